# Update the fixed "Date" footer placeholder text (Insert > Header & Footer >
# Date and time) from 2/6/2023 to 3/20/2023 across every place it is cached:
# the slide master, every slide layout (CustomLayout), and the notes master.

$p = $ppt.ActivePresentation
$newDate = "3/20/2023"
$ppPlaceholderDate = 16

function Update-DateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            if ($shp.TextFrame.TextRange.Text -eq "2/6/2023") {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master
Update-DateShape $p.SlideMaster

# Every slide layout hanging off the slide master
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateShape $layouts.Item($j)
}

# Notes master
Update-DateShape $p.NotesMaster
